# Add a second row of data ("1000.0" in both columns A and B) below the
# existing header row, as plain text values (not numbers), matching the
# original workbook's text-based cell storage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for the numeric-looking strings, then clear the
# resulting explicit number format so no extra style gets attached to
# the cells (keeping them on the default/no style, like the rest of the
# sheet's new content).
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "1000.0"
$ws.Range("A2").Style = "Normal"

$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "1000.0"
$ws.Range("B2").Style = "Normal"
